$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (91 cell edits) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H12").Value = 360.69232
$ws.Range("I12").Value = 234
$ws.Range("K12").Value = 234
$ws.Range("M12").Value = -64
$ws.Range("H17").Value = 1683.6923
$ws.Range("J17").Value = 1683.6923
$ws.Range("L17").Value = 5051.0769
$ws.Range("N17").Value = -5387.0769
$ws.Range("H18").Value = 1471.4286
$ws.Range("I18").Value = 1471.4286
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1471.4286
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -1187.4286
$ws.Range("N18").ClearContents()
$ws.Range("H31").Value = 410
$ws.Range("I31").Value = 410
$ws.Range("K31").Value = 1230
$ws.Range("M31").Value = -1000
$ws.Range("H41").Value = 424
$ws.Range("I41").Value = 473
$ws.Range("K41").Value = 473
$ws.Range("M41").Value = -33
$ws.Range("H53").Value = 1062.4
$ws.Range("I53").Value = 557.25
$ws.Range("J53").Value = 1639.7142
$ws.Range("K53").Value = 557.25
$ws.Range("L53").Value = 1639.7142
$ws.Range("M53").Value = 79.75
$ws.Range("N53").Value = -2913.7142
$ws.Range("H76").Value = 7922.5293
$ws.Range("I76").Value = 9336.125
$ws.Range("K76").Value = 9336.125
$ws.Range("M76").Value = -9021.125
$ws.Range("H79").Value = 7922.5293
$ws.Range("I79").Value = 9336.125
$ws.Range("K79").Value = 9336.125
$ws.Range("M79").Value = -8244.125
$ws.Range("H94").Value = 18666.334
$ws.Range("I94").Value = 15999
$ws.Range("K94").Value = 15999
$ws.Range("M94").Value = -15548
$ws.Range("H95").Value = 114999
$ws.Range("J95").Value = 114999
$ws.Range("L95").Value = 114999
$ws.Range("N95").Value = -120491
$ws.Range("H100").Value = 34212.855
$ws.Range("I100").Value = 16255.286
$ws.Range("K100").Value = 16255.286
$ws.Range("M100").Value = -15714.286
$ws.Range("H115").Value = 691.625
$ws.Range("I115").Value = 630.8333
$ws.Range("K115").Value = 1892.4999
$ws.Range("M115").Value = -325.4999
$ws.Range("H116").Value = 4888.7144
$ws.Range("I116").Value = 4494.773
$ws.Range("J116").Value = 6333.1665
$ws.Range("K116").Value = 4494.773
$ws.Range("L116").Value = 6333.1665
$ws.Range("M116").Value = -1052.773
$ws.Range("N116").Value = -13217.1665
$ws.Range("H118").Value = 313.1111
$ws.Range("J118").Value = 169.66667
$ws.Range("L118").Value = 509.00001
$ws.Range("N118").Value = -3823.00001
$ws.Range("H132").Value = 1511.8438
$ws.Range("I132").Value = 915.25806
$ws.Range("J132").Value = 20006
$ws.Range("K132").Value = 2745.77418
$ws.Range("L132").Value = 60018
$ws.Range("M132").Value = -215.7741799999999
$ws.Range("N132").Value = -65078
$ws.Range("H137").Value = 4654.1665
$ws.Range("I137").Value = 2165.0476
$ws.Range("J137").Value = 10462.111
$ws.Range("K137").Value = 6495.1428
$ws.Range("L137").Value = 31386.333
$ws.Range("M137").Value = -3945.1428
$ws.Range("N137").Value = -36486.333
$ws.Range("H138").Value = 2698.3562
$ws.Range("J138").Value = 2790.8643
$ws.Range("L138").Value = 8372.5929
$ws.Range("N138").Value = -18652.5929
$ws.Range("H141").Value = 6427.857
$ws.Range("I141").Value = 5908.636
$ws.Range("K141").Value = 17725.908
$ws.Range("M141").Value = -12545.908

# ---- Sheet: ARM (63 cell edits) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1437
$ws.Range("I2").Value = 1400.375
$ws.Range("J2").Value = 1632.3334
$ws.Range("K2").Value = 1400.375
$ws.Range("L2").Value = 1632.3334
$ws.Range("M2").Value = -1287.375
$ws.Range("N2").Value = -1858.3334
$ws.Range("H5").Value = 188.2
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 188.2
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 188.2
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -412.2
$ws.Range("H21").Value = 496.33334
$ws.Range("J21").Value = 442.66666
$ws.Range("L21").Value = 442.66666
$ws.Range("N21").Value = -1190.66666
$ws.Range("H32").Value = 7697825.5
$ws.Range("I32").Value = 9618036
$ws.Range("J32").Value = 16982.77
$ws.Range("K32").Value = 9618036
$ws.Range("L32").Value = 16982.77
$ws.Range("M32").Value = -9617749
$ws.Range("N32").Value = -17556.77
$ws.Range("H45").Value = 22729476
$ws.Range("I45").Value = 33335310
$ws.Range("K45").Value = 33335310
$ws.Range("M45").Value = -33334933
$ws.Range("H60").Value = 45500
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H74").Value = 10426496
$ws.Range("I74").Value = 20836136
$ws.Range("K74").Value = 20836136
$ws.Range("M74").Value = -20835262
$ws.Range("H77").Value = 10426496
$ws.Range("I77").Value = 20836136
$ws.Range("K77").Value = 104180680
$ws.Range("M77").Value = -104176312
$ws.Range("H94").Value = 48000
$ws.Range("J94").Value = 48000
$ws.Range("L94").Value = 48000
$ws.Range("N94").Value = -49802
$ws.Range("H110").Value = 1153.3182
$ws.Range("I110").Value = 955.2778
$ws.Range("J110").Value = 2044.5
$ws.Range("K110").Value = 955.2778
$ws.Range("L110").Value = 2044.5
$ws.Range("M110").Value = 1089.7222
$ws.Range("N110").Value = -6134.5
$ws.Range("H116").Value = 1437
$ws.Range("I116").Value = 1400.375
$ws.Range("J116").Value = 1632.3334
$ws.Range("K116").Value = 1400.375
$ws.Range("L116").Value = 1632.3334
$ws.Range("M116").Value = 893.625
$ws.Range("N116").Value = -6220.3334
$ws.Range("H122").Value = 2630.513
$ws.Range("I122").Value = 1767.5834
$ws.Range("K122").Value = 5302.7502
$ws.Range("M122").Value = -2852.7502

# ---- Sheet: BSM (53 cell edits) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1437
$ws.Range("I3").Value = 1400.375
$ws.Range("J3").Value = 1632.3334
$ws.Range("K3").Value = 1400.375
$ws.Range("L3").Value = 1632.3334
$ws.Range("M3").Value = -1286.375
$ws.Range("N3").Value = -1860.3334
$ws.Range("H4").Value = 188.2
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 188.2
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 188.2
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -418.2
$ws.Range("H75").Value = 46657.6
$ws.Range("I75").Value = 6699.5
$ws.Range("J75").Value = 73296.336
$ws.Range("K75").Value = 6699.5
$ws.Range("L75").Value = 73296.336
$ws.Range("M75").Value = -5763.5
$ws.Range("N75").Value = -75168.336
$ws.Range("H78").Value = 46657.6
$ws.Range("I78").Value = 6699.5
$ws.Range("J78").Value = 73296.336
$ws.Range("K78").Value = 20098.5
$ws.Range("L78").Value = 219889.008
$ws.Range("M78").Value = -15418.5
$ws.Range("N78").Value = -229249.008
$ws.Range("H80").Value = 3039.6
$ws.Range("I80").Value = 1733
$ws.Range("J80").Value = 4999.5
$ws.Range("K80").Value = 1733
$ws.Range("L80").Value = 4999.5
$ws.Range("M80").Value = -735
$ws.Range("N80").Value = -6995.5
$ws.Range("H83").Value = 3039.6
$ws.Range("I83").Value = 1733
$ws.Range("J83").Value = 4999.5
$ws.Range("K83").Value = 8665
$ws.Range("L83").Value = 24997.5
$ws.Range("M83").Value = -3673
$ws.Range("N83").Value = -34981.5
$ws.Range("H105").Value = 2062.575
$ws.Range("I105").Value = 1333.55
$ws.Range("J105").Value = 2791.6
$ws.Range("K105").Value = 1333.55
$ws.Range("L105").Value = 2791.6
$ws.Range("M105").Value = 413.45
$ws.Range("N105").Value = -6285.6
$ws.Range("H108").Value = 108999
$ws.Range("J108").Value = 108999
$ws.Range("L108").Value = 108999
$ws.Range("N108").Value = -116679

# ---- Sheet: CRP (113 cell edits) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 233.66667
$ws.Range("I7").Value = 211.90909
$ws.Range("J7").Value = 267.85715
$ws.Range("K7").Value = 211.90909
$ws.Range("L7").Value = 267.85715
$ws.Range("M7").Value = -98.90908999999999
$ws.Range("N7").Value = -493.85715
$ws.Range("H22").Value = 1247
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 994
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 994
$ws.Range("M22").Value = -1150
$ws.Range("N22").Value = -1694
$ws.Range("H31").Value = 528490.6
$ws.Range("I31").Value = 11418.934
$ws.Range("J31").Value = 787026.4399999999
$ws.Range("K31").Value = 11418.934
$ws.Range("L31").Value = 787026.4399999999
$ws.Range("M31").Value = -11123.934
$ws.Range("N31").Value = -787616.4399999999
$ws.Range("H34").Value = 528490.6
$ws.Range("I34").Value = 11418.934
$ws.Range("J34").Value = 787026.4399999999
$ws.Range("K34").Value = 11418.934
$ws.Range("L34").Value = 787026.4399999999
$ws.Range("M34").Value = -11216.934
$ws.Range("N34").Value = -787430.4399999999
$ws.Range("H76").Value = 8721.4
$ws.Range("I76").Value = 8721.4
$ws.Range("K76").Value = 8721.4
$ws.Range("M76").Value = -8406.4
$ws.Range("H79").Value = 8721.4
$ws.Range("I79").Value = 8721.4
$ws.Range("K79").Value = 8721.4
$ws.Range("M79").Value = -7629.4
$ws.Range("H86").Value = 149753.58
$ws.Range("I86").Value = 8249.5
$ws.Range("J86").Value = 206355.2
$ws.Range("K86").Value = 8249.5
$ws.Range("L86").Value = 206355.2
$ws.Range("M86").Value = -7126.5
$ws.Range("N86").Value = -208601.2
$ws.Range("H89").Value = 149753.58
$ws.Range("I89").Value = 8249.5
$ws.Range("J89").Value = 206355.2
$ws.Range("K89").Value = 41247.5
$ws.Range("L89").Value = 1031776
$ws.Range("M89").Value = -35631.5
$ws.Range("N89").Value = -1043008
$ws.Range("H94").Value = 3361.8
$ws.Range("I94").Value = 2816.8
$ws.Range("J94").Value = 3543.4666
$ws.Range("K94").Value = 2816.8
$ws.Range("L94").Value = 3543.4666
$ws.Range("M94").Value = -2365.8
$ws.Range("N94").Value = -4445.4666
$ws.Range("H99").Value = 4011.5
$ws.Range("I99").Value = 4011.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4011.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2513.5
$ws.Range("N99").ClearContents()
$ws.Range("H103").Value = 53497.5
$ws.Range("J103").Value = 99995
$ws.Range("L103").Value = 99995
$ws.Range("N103").Value = -102339
$ws.Range("H107").Value = 2525.7222
$ws.Range("I107").Value = 959.1539
$ws.Range("J107").Value = 6598.8
$ws.Range("K107").Value = 959.1539
$ws.Range("L107").Value = 6598.8
$ws.Range("M107").Value = 960.8461
$ws.Range("N107").Value = -10438.8
$ws.Range("H112").Value = 111409.664
$ws.Range("J112").Value = 111409.664
$ws.Range("L112").Value = 111409.664
$ws.Range("N112").Value = -114363.664
$ws.Range("H114").Value = 9500
$ws.Range("J114").Value = 9500
$ws.Range("L114").Value = 9500
$ws.Range("N114").Value = -18178
$ws.Range("H116").Value = 90577.336
$ws.Range("J116").Value = 90577.336
$ws.Range("L116").Value = 90577.336
$ws.Range("N116").Value = -99755.336
$ws.Range("H119").Value = 52160.332
$ws.Range("J119").Value = 52160.332
$ws.Range("L119").Value = 52160.332
$ws.Range("N119").Value = -61836.332
$ws.Range("H124").Value = 55941.145
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("H126").Value = 4011.5
$ws.Range("I126").Value = 4011.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 12034.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9564.5
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 10089.667
$ws.Range("J132").Value = 22513.25
$ws.Range("L132").Value = 67539.75
$ws.Range("N132").Value = -72599.75
$ws.Range("H134").Value = 5298.4287
$ws.Range("I134").Value = 2022.75
$ws.Range("J134").Value = 9666
$ws.Range("K134").Value = 6068.25
$ws.Range("L134").Value = 28998
$ws.Range("M134").Value = -3533.25
$ws.Range("N134").Value = -34068

# ---- Sheet: CUL (34 cell edits) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1075.2222
$ws.Range("J23").Value = 695.4
$ws.Range("L23").Value = 2086.2
$ws.Range("N23").Value = -2556.2
$ws.Range("H39").Value = 141
$ws.Range("I39").Value = 141
$ws.Range("K39").Value = 423
$ws.Range("M39").Value = -129
$ws.Range("H51").Value = 36813.75
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 36813.75
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 110441.25
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -111361.25
$ws.Range("H80").Value = 4527.077
$ws.Range("J80").Value = 4531.5
$ws.Range("L80").Value = 13594.5
$ws.Range("N80").Value = -15466.5
$ws.Range("H83").Value = 4527.077
$ws.Range("J83").Value = 4531.5
$ws.Range("L83").Value = 40783.5
$ws.Range("N83").Value = -50143.5
$ws.Range("H113").Value = 1332.3684
$ws.Range("I113").Value = 739.8
$ws.Range("J113").Value = 1544
$ws.Range("K113").Value = 2219.4
$ws.Range("L113").Value = 4632
$ws.Range("M113").Value = -49.39999999999964
$ws.Range("N113").Value = -8972
$ws.Range("H131").Value = 18775.834
$ws.Range("J131").Value = 20410
$ws.Range("L131").Value = 61230
$ws.Range("N131").Value = -71310

# ---- Sheet: GSM (30 cell edits) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 55831.168
$ws.Range("J93").Value = 55831.168
$ws.Range("L93").Value = 55831.168
$ws.Range("N93").Value = -59575.168
$ws.Range("H95").Value = 40571
$ws.Range("J95").Value = 40571
$ws.Range("L95").Value = 40571
$ws.Range("N95").Value = -46063
$ws.Range("H99").Value = 53998
$ws.Range("I99").Value = 1996
$ws.Range("J99").Value = 106000
$ws.Range("K99").Value = 1996
$ws.Range("L99").Value = 106000
$ws.Range("M99").Value = 250
$ws.Range("N99").Value = -110492
$ws.Range("H102").Value = 1754.8684
$ws.Range("I102").Value = 909.7143
$ws.Range("K102").Value = 909.7143
$ws.Range("M102").Value = 712.2857
$ws.Range("H113").Value = 3741.7646
$ws.Range("I113").Value = 2997.75
$ws.Range("K113").Value = 2997.75
$ws.Range("M113").Value = -827.75
$ws.Range("H132").Value = 41669532
$ws.Range("I132").Value = 52634410
$ws.Range("J132").Value = 3001.2
$ws.Range("K132").Value = 157903230
$ws.Range("L132").Value = 9003.599999999999
$ws.Range("M132").Value = -157900700
$ws.Range("N132").Value = -14063.6

# ---- Sheet: LTW (46 cell edits) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H7").Value = 12756.85
$ws.Range("I7").Value = 11673.625
$ws.Range("J7").Value = 13479
$ws.Range("K7").Value = 11673.625
$ws.Range("L7").Value = 13479
$ws.Range("M7").Value = -11561.625
$ws.Range("N7").Value = -13703
$ws.Range("H40").Value = 4044.3447
$ws.Range("I40").Value = 2538.2727
$ws.Range("J40").Value = 8777.714
$ws.Range("K40").Value = 2538.2727
$ws.Range("L40").Value = 8777.714
$ws.Range("M40").Value = -2402.2727
$ws.Range("N40").Value = -9049.714
$ws.Range("H61").Value = 1026.3043
$ws.Range("I61").Value = 645.1875
$ws.Range("J61").Value = 1897.4286
$ws.Range("K61").Value = 645.1875
$ws.Range("L61").Value = 1897.4286
$ws.Range("M61").Value = -443.1875
$ws.Range("N61").Value = -2301.4286
$ws.Range("H113").Value = 1026.3043
$ws.Range("I113").Value = 645.1875
$ws.Range("J113").Value = 1897.4286
$ws.Range("K113").Value = 645.1875
$ws.Range("L113").Value = 1897.4286
$ws.Range("M113").Value = 1524.8125
$ws.Range("N113").Value = -6237.4286
$ws.Range("H122").Value = 5425.61
$ws.Range("I122").Value = 4757.375
$ws.Range("J122").Value = 6369
$ws.Range("K122").Value = 14272.125
$ws.Range("L122").Value = 19107
$ws.Range("M122").Value = -11822.125
$ws.Range("N122").Value = -24007
$ws.Range("H126").Value = 12756.85
$ws.Range("I126").Value = 11673.625
$ws.Range("J126").Value = 13479
$ws.Range("K126").Value = 35020.875
$ws.Range("L126").Value = 40437
$ws.Range("M126").Value = -32550.875
$ws.Range("N126").Value = -45377

# ---- Sheet: WVR (45 cell edits) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 55177.5
$ws.Range("I6").Value = 110005
$ws.Range("J6").Value = 350
$ws.Range("K6").Value = 110005
$ws.Range("L6").Value = 350
$ws.Range("M6").Value = -109890
$ws.Range("N6").Value = -580
$ws.Range("H11").Value = 34943.5
$ws.Range("J11").Value = 34943.5
$ws.Range("L11").Value = 34943.5
$ws.Range("N11").Value = -35227.5
$ws.Range("H15").Value = 36669
$ws.Range("J15").Value = 36669
$ws.Range("L15").Value = 36669
$ws.Range("N15").Value = -37245
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H107").Value = 31251564
$ws.Range("I107").Value = 45456040
$ws.Range("J107").Value = 1712.4
$ws.Range("K107").Value = 136368120
$ws.Range("L107").Value = 5137.200000000001
$ws.Range("M107").Value = -136366200
$ws.Range("N107").Value = -8977.200000000001
$ws.Range("H113").Value = 970.86664
$ws.Range("I113").Value = 825.5
$ws.Range("J113").Value = 1261.6
$ws.Range("K113").Value = 2476.5
$ws.Range("L113").Value = 3784.8
$ws.Range("M113").Value = -306.5
$ws.Range("N113").Value = -8124.799999999999
$ws.Range("H122").Value = 3675.2666
$ws.Range("I122").Value = 3723.5
$ws.Range("K122").Value = 11170.5
$ws.Range("M122").Value = -8720.5
$ws.Range("H126").Value = 2850
$ws.Range("I126").Value = 2400
$ws.Range("K126").Value = 7200
$ws.Range("M126").Value = -4730
$ws.Range("H132").Value = 259988.89
$ws.Range("I132").Value = 2825.5
$ws.Range("K132").Value = 8476.5
$ws.Range("M132").Value = -5946.5
